# Render the automatic reports
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new "NA" entries appear on row 3 (auto-generated report cells)
$ws.Range("W3").Value = "NA"
$ws.Range("AN3").Value = "NA"

# Scroll the view one column further right and move the selection
$excel.ActiveWindow.ScrollColumn = 15
$ws.Range("AM4").Select()
